# Restored cell C10 ("R30" rule, "Integer min" / From-hour) on the
# "Rules" sheet from 18 back to 1, per the commit's target revision.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Range("C10").Value = 1
